$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.726.42"
$ws.Range("E2").Value = "  +0.87%  "
$ws.Range("D3").Value = "2.100.57"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("E5").Value = "  -0.19%  "
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.97"
$ws.Range("E7").Value = "  +1.11%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +2.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0844"
$ws.Range("E10").Value = "  -0.13%  "
$ws.Range("E11").Value = "  -0.44%  "
$ws.Range("E12").Value = "  +7.06%  "
$ws.Range("D13").Value = "2.411.67"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.06"
$ws.Range("E14").Value = "  -1.12%  "
$ws.Range("E15").Value = "  +3.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.52"
$ws.Range("E16").Value = "  +1.13%  "
$ws.Range("D17").Value = "2.092.42"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").Value = "38.719.48"
$ws.Range("E18").Value = "  +1.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.84"
$ws.Range("E19").Value = "  +2.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.08"
$ws.Range("E20").Value = "  +1.14%  "
$ws.Range("E21").Value = "  +0.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.50"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.37"
$ws.Range("E24").Value = "  -2.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.31"
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "172.05"
$ws.Range("E26").Value = "  +1.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.56"
$ws.Range("E27").Value = "  +1.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.138"
$ws.Range("E28").Value = "  +6.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.43"
$ws.Range("E29").Value = "  +4.89%  "
$ws.Range("E30").Value = "  +1.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.55"
$ws.Range("E31").Value = "  +8.14%  "
$ws.Range("E32").Value = "  +0.68%  "
$ws.Range("E33").Value = "  +2.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.76"
$ws.Range("E34").Value = "  +0.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.96"
$ws.Range("E35").Value = "  +8.54%  "
$ws.Range("E36").Value = "  +2.31%  "
$ws.Range("E37").Value = "  +1.75%  "
$ws.Range("E38").Value = "  +3.20%  "
$ws.Range("E39").Value = "  -0.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.11"
$ws.Range("E40").Value = "  -0.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "102.77"
$ws.Range("E41").Value = "  +2.97%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0228"
$ws.Range("E42").Value = "  +4.26%  "
$ws.Range("D43").Value = "1.535.15"
$ws.Range("E43").Value = "  -0.15%  "
$ws.Range("E44").Value = "  +4.67%  "
$ws.Range("E45").Value = "  -1.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.16"
$ws.Range("E46").Value = "  +3.81%  "
$ws.Range("E47").Value = "  +0.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.12"
$ws.Range("E48").Value = "  -0.94%  "
$ws.Range("E49").Value = "  +1.02%  "
$ws.Range("E50").Value = "  -0.38%  "
$ws.Range("D51").Value = "2.298.01"
$ws.Range("E51").Value = "  +0.23%  "
